$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting rows 59:79 down to 60:80
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly record
$ws.Cells.Item(59, 1).Value = 10
$ws.Cells.Item(59, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(59, 3).Value = "La Araucanía"
$ws.Cells.Item(59, 4).Value = 44809
$ws.Cells.Item(59, 5).Value = 9
$ws.Cells.Item(59, 6).Value = 300000001
$ws.Cells.Item(59, 7).Value = "Rabanito"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 30
$ws.Cells.Item(59, 11).Value = 10000
$ws.Cells.Item(59, 12).Value = 10000
$ws.Cells.Item(59, 13).Value = 10000
$ws.Cells.Item(59, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(59, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(59, 16).Value = 833
$ws.Cells.Item(59, 17).Value = 12
$ws.Cells.Item(59, 18).Value = "Hortaliza"
